$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 9 (un_franzosa_ControlvsCD_Fp),
# shifting it and everything below down by one.
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the un_franzosa_ControlvsCD_ConvCD data.
$ws.Range("A9").Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.4
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6
$ws.Range("H9").Value = 0.6

# After the first insert, un_franzosa_ControlvsUC_Age is at row 14 and
# un_franzosa_ControlvsUC_Fp is at row 15. Insert a new row before the Fp row.
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the un_franzosa_ControlvsUC_ConvUC data.
$ws.Range("A15").Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.4
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6
$ws.Range("H15").Value = 0.6
